$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.66138833333333
$ws.Range("H2").Value = 64.98416499999999
$ws.Range("I2").Value = 0.2252765553546639
$ws.Range("J2").Value = 0.2252765553546639
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.535604
$ws.Range("N2").Value = 49.606812
$ws.Range("O2").Value = 0.2120453146491552
$ws.Range("P2").Value = 0.2120453146491552
$ws.Range("Q2").Value = 358.1841395702199
$ws.Range("R2").Value = 3223.657256131979
$ws.Range("S2").Value = 0.04776883806325753
$ws.Range("T2").Value = 0.04776883806325755

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.66138833333333
$ws.Range("H3").Value = 64.98416499999999
$ws.Range("I3").Value = 0.2252765553546639
$ws.Range("J3").Value = 0.2252765553546639
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 40.62063066666667
$ws.Range("N3").Value = 121.861892
$ws.Range("O3").Value = 0.5209011059384622
$ws.Range("P3").Value = 0.5209011059384622
$ws.Range("Q3").Value = 879.8992552155755
$ws.Range("R3").Value = 7919.093296940179
$ws.Range("S3").Value = 0.1173468068262516
$ws.Range("T3").Value = 0.1173468068262516

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.66138833333333
$ws.Range("H4").Value = 64.98416499999999
$ws.Range("I4").Value = 0.2252765553546639
$ws.Range("J4").Value = 0.2252765553546639
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 20.825229
$ws.Range("N4").Value = 62.475687
$ws.Range("O4").Value = 0.2670535794123827
$ws.Range("P4").Value = 0.2670535794123827
$ws.Range("Q4").Value = 451.103372499595
$ws.Range("R4").Value = 4059.930352496354
$ws.Range("S4").Value = 0.06016091046515475
$ws.Range("T4").Value = 0.06016091046515475

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 24.68088566666666
$ws.Range("H5").Value = 74.04265699999999
$ws.Range("I5").Value = 0.2566790650963491
$ws.Range("J5").Value = 0.2566790650963491
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 16.535604
$ws.Range("N5").Value = 49.606812
$ws.Range("O5").Value = 0.2120453146491552
$ws.Range("P5").Value = 0.2120453146491552
$ws.Range("Q5").Value = 408.1133517532759
$ws.Range("R5").Value = 3673.020165779483
$ws.Range("S5").Value = 0.05442759312220633
$ws.Range("T5").Value = 0.05442759312220635

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 24.68088566666666
$ws.Range("H6").Value = 74.04265699999999
$ws.Range("I6").Value = 0.2566790650963491
$ws.Range("J6").Value = 0.2566790650963491
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 40.62063066666667
$ws.Range("N6").Value = 121.861892
$ws.Range("O6").Value = 0.5209011059384622
$ws.Range("P6").Value = 0.5209011059384622
$ws.Range("Q6").Value = 1002.553141191894
$ws.Range("R6").Value = 9022.978270727044
$ws.Range("S6").Value = 0.1337044088799388
$ws.Range("T6").Value = 0.1337044088799388

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 24.68088566666666
$ws.Range("H7").Value = 74.04265699999999
$ws.Range("I7").Value = 0.2566790650963491
$ws.Range("J7").Value = 0.2566790650963491
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 20.825229
$ws.Range("N7").Value = 62.475687
$ws.Range("O7").Value = 0.2670535794123827
$ws.Range("P7").Value = 0.2670535794123827
$ws.Range("Q7").Value = 513.985095931151
$ws.Range("R7").Value = 4625.865863380359
$ws.Range("S7").Value = 0.068547063094204
$ws.Range("T7").Value = 0.06854706309420401

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 49.812376
$ws.Range("H8").Value = 149.437128
$ws.Range("I8").Value = 0.518044379548987
$ws.Range("J8").Value = 0.5180443795489871
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 16.535604
$ws.Range("N8").Value = 49.606812
$ws.Range("O8").Value = 0.2120453146491552
$ws.Range("P8").Value = 0.2120453146491552
$ws.Range("Q8").Value = 823.677723835104
$ws.Range("R8").Value = 7413.099514515936
$ws.Range("S8").Value = 0.1098488834636913
$ws.Range("T8").Value = 0.1098488834636914

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 49.812376
$ws.Range("H9").Value = 149.437128
$ws.Range("I9").Value = 0.518044379548987
$ws.Range("J9").Value = 0.5180443795489871
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 40.62063066666667
$ws.Range("N9").Value = 121.861892
$ws.Range("O9").Value = 0.5209011059384622
$ws.Range("P9").Value = 0.5209011059384622
$ws.Range("Q9").Value = 2023.410128125131
$ws.Range("R9").Value = 18210.69115312618
$ws.Range("S9").Value = 0.2698498902322717
$ws.Range("T9").Value = 0.2698498902322718

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 49.812376
$ws.Range("H10").Value = 149.437128
$ws.Range("I10").Value = 0.518044379548987
$ws.Range("J10").Value = 0.5180443795489871
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 20.825229
$ws.Range("N10").Value = 62.475687
$ws.Range("O10").Value = 0.2670535794123827
$ws.Range("P10").Value = 0.2670535794123827
$ws.Range("Q10").Value = 1037.354137234104
$ws.Range("R10").Value = 9336.187235106936
$ws.Range("S10").Value = 0.1383456058530239
$ws.Range("T10").Value = 0.1383456058530239
